# Fix German encoding issues (mojibake) that appeared in the Excel export.
# Some umlaut / sharp-s characters (e.g. "ü", "ß") were lost in a previous
# export and replaced with the Unicode replacement character (U+FFFD).
# This function restores the most common German special characters based on
# their surrounding context.
function Fix-GermanEncoding($Text) {
    if ([string]::IsNullOrEmpty($Text)) {
        return $Text
    }

    $replacementChar = [char]0xFFFD
    $ue = [char]0xFC   # ü
    $ss = [char]0xDF   # ß

    $fixed = $Text

    $pattern1 = $replacementChar + "ringe"
    $value1 = $ue + "ringe"
    $fixed = $fixed -replace $pattern1, $value1               # Th?ringen -> Thüringen

    $pattern2 = "wei" + $replacementChar + "bach"
    $value2 = "wei" + $ss + "bach"
    $fixed = $fixed -replace $pattern2, $value2                # wei?bach -> weißbach

    $pattern3 = "stra" + $replacementChar + "e"
    $value3 = "stra" + $ss + "e"
    $fixed = $fixed -replace $pattern3, $value3                 # Nordstra?e -> Nordstraße

    $pattern4 = "He" + $replacementChar + " "
    $value4 = "He" + $ss + " "
    $fixed = $fixed -replace $pattern4, $value4                  # He? -> Heß

    return $fixed
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$replacementChar = [char]0xFFFD

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $value = $cell.Value2
        if (($value -is [string]) -and ($value.IndexOf($replacementChar) -ge 0)) {
            $newValue = Fix-GermanEncoding $value
            $cell.Value2 = $newValue
        }
    }
}

# Update the Excel export: column K (which holds the longest corrected
# strings) needs its width refreshed to fit the updated content.
$ws.Columns.Item(11).ColumnWidth = 92.8

$wb.Save()
